$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns keep their original text format
# (values such as "1.009" or "29.465.74" must not be auto-converted to numbers)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "29.465.74"
$ws.Range("E2").Value = "  -3.11%  "
$ws.Range("D3").Value = "1.992.96"
$ws.Range("E3").Value = "  -6.28%  "
$ws.Range("D4").Value = "1.009"
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "329.41"
$ws.Range("E5").Value = "  -5.27%  "
$ws.Range("D6").Value = "1.007"
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("D7").Value = "0.5014"
$ws.Range("E7").Value = "  -3.94%  "
$ws.Range("D8").Value = "0.4222"
$ws.Range("E8").Value = "  -5.84%  "
$ws.Range("D9").Value = "52.35"
$ws.Range("E9").Value = "  -3.00%  "
$ws.Range("D10").Value = "0.08903"
$ws.Range("E10").Value = "  -5.02%  "
$ws.Range("E11").Value = "  -5.45%  "
$ws.Range("D12").Value = "23.33"
$ws.Range("E12").Value = "  -8.49%  "
$ws.Range("D13").Value = "8.132"
$ws.Range("E13").Value = "  -6.29%  "
$ws.Range("D14").Value = "1.989.39"
$ws.Range("E14").Value = "  -6.45%  "
$ws.Range("D15").Value = "6.514"
$ws.Range("E15").Value = "  -6.70%  "
$ws.Range("D16").Value = "96.26"
$ws.Range("E16").Value = "  -6.67%  "
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("D18").Value = "0.00001108"
$ws.Range("E18").Value = "  -5.76%  "
$ws.Range("D19").Value = "0.06621"
$ws.Range("E19").Value = "  -1.24%  "
$ws.Range("E20").Value = "  -8.55%  "
$ws.Range("D21").Value = "1.006"
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").Value = "5.976"
$ws.Range("E22").Value = "  -5.46%  "
$ws.Range("D23").Value = "29.495.51"
$ws.Range("E23").Value = "  -3.02%  "
$ws.Range("D24").Value = "11.91"
$ws.Range("E24").Value = "  -6.91%  "
$ws.Range("D25").Value = "2.274"
$ws.Range("E25").Value = "  -2.43%  "
$ws.Range("D26").Value = "157.73"
$ws.Range("E26").Value = "  -3.49%  "
$ws.Range("E27").Value = "  -7.35%  "
$ws.Range("D28").Value = "6.539"
$ws.Range("E28").Value = "  -4.47%  "
$ws.Range("D29").Value = "2.343"
$ws.Range("E29").Value = "  -8.38%  "
$ws.Range("D30").Value = "128.00"
$ws.Range("E30").Value = "  -4.83%  "
$ws.Range("D31").Value = "1.055"
$ws.Range("E31").Value = "  -9.75%  "
$ws.Range("D32").Value = "0.09960"
$ws.Range("E32").Value = "  -6.04%  "
$ws.Range("D33").Value = "1.558"
$ws.Range("E33").Value = "  -13.78%  "
$ws.Range("D34").Value = "5.848"
$ws.Range("E34").Value = "  -7.38%  "
$ws.Range("D35").Value = "3.785"
$ws.Range("E35").Value = "  -4.32%  "
$ws.Range("D36").Value = "9.606"
$ws.Range("E36").Value = "  -10.56%  "
$ws.Range("D37").Value = "0.02463"
$ws.Range("E37").Value = "  -6.96%  "
$ws.Range("E38").Value = "  -7.78%  "
$ws.Range("E39").Value = "  -3.97%  "
$ws.Range("D40").Value = "0.6541"
$ws.Range("E40").Value = "  -8.77%  "
$ws.Range("D41").Value = "11.75"
$ws.Range("E41").Value = "  -8.31%  "
$ws.Range("E42").Value = "  -8.40%  "
$ws.Range("D43").Value = "1.005"
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("D44").Value = "0.6357"
$ws.Range("E44").Value = "  -8.92%  "
$ws.Range("D45").Value = "2.222"
$ws.Range("E45").Value = "  -7.47%  "
$ws.Range("D46").Value = "13.41"
$ws.Range("E46").Value = "  -9.28%  "
$ws.Range("E47").Value = "  -0.14%  "
$ws.Range("D48").Value = "3.517"
$ws.Range("E48").Value = "  -3.20%  "
$ws.Range("D49").Value = "0.00000000338"
$ws.Range("E49").Value = "  -3.01%  "
$ws.Range("D50").Value = "0.07001"
$ws.Range("E50").Value = "  -3.15%  "
$ws.Range("E51").Value = "  -6.11%  "
